$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Sheet index 2: Restricciones_del_lider
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("A2") "2.3000000000000003 - x"
Set-TextValue $ws2.Range("B2") "-3.3000000000000003"
Set-TextValue $ws2.Range("D2") "0.51"
Set-TextValue $ws2.Range("A3") "-2.3000000000000003 + x"
Set-TextValue $ws2.Range("B3") "1.3000000000000003"
Set-TextValue $ws2.Range("D3") "0.17"

# Sheet index 3: Restricciones_del_follower
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("A2") "-10.877777777777776 + 2.4444444444444446y"
Set-TextValue $ws3.Range("B2") "9.877777777777776"
Set-TextValue $ws3.Range("D2") "0.82"
Set-TextValue $ws3.Range("E2") "3.7"
Set-TextValue $ws3.Range("F2") "5.5"
Set-TextValue $ws3.Range("A3") "-5.562499999999999 + 1.25y"
Set-TextValue $ws3.Range("B3") "4.562499999999999"
Set-TextValue $ws3.Range("D3") "0.81"
Set-TextValue $ws3.Range("E3") "7.199999999999999"
Set-TextValue $ws3.Range("F3") "1.0"

# Sheet index 4: Punto_modificado
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "2.3000000000000003"
Set-TextValue $ws4.Range("B2") "4.449999999999999"

# Sheet index 5: Vector_bf
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-1.084944444444441"

# Sheet index 6: Vector_BF
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "25.035000000000007"
Set-TextValue $ws6.Range("A3") "-22.219444444444445"

# Sheet index 7: Vector_Alpha
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.25
